# Update the "Förändrad" (Changed) date column from 2024-02-09 (45331)
# to 2024-02-10 (45332) for all data rows (rows 2-27, column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 27 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45331) {
        $cell.Value2 = 45332
    }
}
